# updating new version of table2
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Re-label header row (typed D,C,B,A so the shared-string table ends
#    up ordered Age, Pulse, SugarBlood, ID) and drop the old
#    "illBefore" column (E).
# ---------------------------------------------------------------------
$ws.Range("D1").Value = "Age"
$ws.Range("C1").Value = "Pulse"
$ws.Range("B1").Value = "SugarBlood"
$ws.Range("A1").Value = "ID"

$ws.Range("E1:E11").Clear()

# ---------------------------------------------------------------------
# 2) Make every data row share the same base look (the style used by
#    row 2) before editing values, so rows that used to be
#    "highlighted" (rows 5,6,9,10) become plain like the rest.
# ---------------------------------------------------------------------
$ws.Range("A2:D2").Copy()
$ws.Range("A5:D6").PasteSpecial(-4122)
$ws.Range("A9:D10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3) Update the Pulse (column C) values that changed.
# ---------------------------------------------------------------------
$ws.Range("C4").Value = 88
$ws.Range("C6").Value = 83
$ws.Range("C7").Value = 82
$ws.Range("C8").Value = 81
$ws.Range("C10").Value = 89
$ws.Range("C11").Value = 69

# ---------------------------------------------------------------------
# 4) Append five new patient rows (12-16), reusing row 11's formatting.
# ---------------------------------------------------------------------
$ws.Range("A11:D11").Copy()
$ws.Range("A12:D16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$newRows = @(
  @(11,158,71,24),
  @(12,88,73,25),
  @(13,89,77,19),
  @(14,68,78,18),
  @(15,72,87,51)
)
$r = 12
foreach ($row in $newRows) {
  $ws.Cells.Item($r,1).Value = $row[0]
  $ws.Cells.Item($r,2).Value = $row[1]
  $ws.Cells.Item($r,3).Value = $row[2]
  $ws.Cells.Item($r,4).Value = $row[3]
  $r = $r + 1
}

# ---------------------------------------------------------------------
# 5) Work out the two final cell looks (plain/centered data cells and
#    bold/Lucida Handwriting/centered header cells) on scratch cells
#    first, then fan them out with a format-only paste. This keeps the
#    number of newly generated font/style records to a minimum instead
#    of forking a brand new style for every individual cell touched.
# ---------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$dataLook = $ws.Range("Z1")
$dataLook.HorizontalAlignment = -4108

$ws.Range("A1").Copy()
$ws.Range("Z2").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$headerLook = $ws.Range("Z2")
$headerLook.Font.Bold = $true
$headerLook.Font.Name = "Lucida Handwriting"
$headerLook.HorizontalAlignment = -4108

$dataLook.Copy()
$ws.Range("A2:D16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$headerLook.Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("Z1:Z2").Clear()

# Header row is a touch taller in the new layout.
$ws.Rows.Item(1).RowHeight = 15

# ---------------------------------------------------------------------
# 6) Match the saved selection/cursor position.
# ---------------------------------------------------------------------
$ws.Range("C4").Select()
